# Update the fixed "Date Placeholder" text (10/20/25 -> 10/21/25) on the
# slide master and on every slide layout, and update the GitHub repo URL
# on slide 2 from the 2024 seminar to the 2025 seminar.

$p = $ppt.ActivePresentation

$oldDate = "10/20/25"
$newDate = "10/21/25"

function Update-DatePlaceholder($shapes) {
    for ($i = 1; $i -le $shapes.Count; $i++) {
        $shp = $shapes.Item($i)
        if ($shp.HasTextFrame) {
            if ($shp.TextFrame.HasText) {
                if ($shp.TextFrame.TextRange.Text -eq $oldDate) {
                    $shp.TextFrame.TextRange.Text = $newDate
                }
            }
        }
    }
}

# Slide master
$master = $p.SlideMaster
Update-DatePlaceholder $master.Shapes

# Every slide layout belonging to the master
for ($li = 1; $li -le $master.CustomLayouts.Count; $li++) {
    $layout = $master.CustomLayouts.Item($li)
    Update-DatePlaceholder $layout.Shapes
}

# Slide 2: update the GitHub repository URL for the 2025 seminar
$s2 = $p.Slides.Item(2)
for ($i = 1; $i -le $s2.Shapes.Count; $i++) {
    $shp = $s2.Shapes.Item($i)
    if ($shp.HasTextFrame) {
        if ($shp.TextFrame.HasText) {
            $t = $shp.TextFrame.TextRange.Text
            if ($t -eq "https://github.com/Networks-Learning/hcml-seminar-2024") {
                $shp.TextFrame.TextRange.Text = "https://github.com/Networks-Learning/hcml-seminar-2025"
            }
        }
    }
}
